# Update the marksheet's "Marking" (B11) and "Total" (B12) correct-mark
# values, and the corresponding Corr/Total display string in E12.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 95
$ws.Range("E12").Value = "95/140"
